$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows: Índice, Distancia, max, min, Tempo
$data = @(
    @(0, 3435.566666666667, 3512, 3399, 0.03300157388051351),
    @(1, 3234.333333333333, 3387, 3029, 0.03539212544759115),
    @(2, 3458.2, 3585, 3238, 0.03497441609700521),
    @(3, 3570.733333333333, 3917, 3312, 0.03509840965270996),
    @(4, 3298, 3555, 2981, 0.03635515371958415),
    @(5, 4081.466666666667, 4227, 3990, 0.03608304659525553),
    @(6, 3107.2, 3344, 2966, 0.0347800334294637),
    @(7, 3790.366666666667, 4021, 3582, 0.0362070878346761),
    @(8, 3149.266666666667, 3377, 2935, 0.03646350701649984),
    @(9, 3122.9, 3531, 3029, 0.03254377841949463)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
